$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-05-12 Monday"

# Update each equation cell in the table (20 rows x 5 columns, row-major order)
$t = $d.Tables.Item(1)
$values = @(
    "7+31=",
    "11+75=",
    "43+23=",
    "16+54=",
    "25+41=",
    "44-18=",
    "51+24=",
    "54+1=",
    "39-2=",
    "72-14=",
    "63-58=",
    "55-46=",
    "15+81=",
    "42+31=",
    "79-18=",
    "63+12=",
    "13+38=",
    "23+31=",
    "38+38=",
    "31+0=",
    "30-9=",
    "53-7=",
    "71-17=",
    "15-8=",
    "52-9=",
    "77-72=",
    "1+39=",
    "43+6=",
    "35-4=",
    "60-8=",
    "39-36=",
    "36-21=",
    "18+18=",
    "37+51=",
    "59-17=",
    "69+27=",
    "57+15=",
    "42-40=",
    "10+40=",
    "29-29=",
    "54-33=",
    "89-3=",
    "86-43=",
    "62-10=",
    "89-0=",
    "25+32=",
    "17+37=",
    "54+6=",
    "10+76=",
    "73+26=",
    "67-63=",
    "46+25=",
    "85-30=",
    "0+29=",
    "54+32=",
    "35-16=",
    "19+30=",
    "33+9=",
    "58-26=",
    "50+23=",
    "7+7=",
    "23+10=",
    "96-30=",
    "66+10=",
    "51+10=",
    "50-31=",
    "6+13=",
    "77-38=",
    "20-10=",
    "49-22=",
    "0+0=",
    "57-14=",
    "3+71=",
    "16+29=",
    "78-31=",
    "63-9=",
    "69-14=",
    "35-31=",
    "80-34=",
    "67-47=",
    "96-59=",
    "87-46=",
    "92-70=",
    "9+29=",
    "42+42=",
    "24+63=",
    "10+31=",
    "34-9=",
    "80-39=",
    "44+45=",
    "47-21=",
    "4+4=",
    "20+78=",
    "54-34=",
    "96-69=",
    "62-14=",
    "43+2=",
    "74-73=",
    "51+31=",
    "46-40="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output ("Updated " + $idx + " cells")